$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") date values for rows 2-180 are all updated
# from serial date 45178 (2023-09-09) to 45179 (2023-09-10).
$ws.Range("C2:C180").Value = 45179
